$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 49, shifting existing rows 49-80 down to 51-82
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Insert()

# Row 49 - new data
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 45089
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100107
$ws.Range("H49").Value = "Otros"
$ws.Range("I49").Value = 100107001
$ws.Range("J49").Value = "Caqui"
$ws.Range("K49").Value = "Fuyu"
$ws.Range("L49").Value = "Especial"
$ws.Range("M49").Value = 450
$ws.Range("N49").Value = 13000
$ws.Range("O49").Value = 16000
$ws.Range("P49").Value = 14667
$ws.Range("Q49").Value = "`$/caja 10 kilos granel"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 14667
$ws.Range("T49").Value = 1

# Row 50 - new data
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 45089
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100107
$ws.Range("H50").Value = "Otros"
$ws.Range("I50").Value = 100107001
$ws.Range("J50").Value = "Caqui"
$ws.Range("K50").Value = "Fuyu"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 400
$ws.Range("N50").Value = 10000
$ws.Range("O50").Value = 10000
$ws.Range("P50").Value = 10000
$ws.Range("Q50").Value = "`$/caja 10 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 10000
$ws.Range("T50").Value = 1
